# Update the time_taken timestamps on the "data" sheet to reflect the
# re-run export time.
$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("data")

$newTimestamps = @(
    "2021-10-05 14:19:16.021064",
    "2021-10-05 14:19:16.021073",
    "2021-10-05 14:19:16.021076",
    "2021-10-05 14:19:16.021079",
    "2021-10-05 14:19:16.021082",
    "2021-10-05 14:19:16.021085",
    "2021-10-05 14:19:16.021088",
    "2021-10-05 14:19:16.021091",
    "2021-10-05 14:19:16.021094",
    "2021-10-05 14:19:16.021096",
    "2021-10-05 14:19:16.021099",
    "2021-10-05 14:19:16.021102",
    "2021-10-05 14:19:16.021105",
    "2021-10-05 14:19:16.021107",
    "2021-10-05 14:19:16.021110"
)

for ($i = 0; $i -lt $newTimestamps.Length; $i++) {
    $row = $i + 2
    $wsData.Cells.Item($row, 6).Value = $newTimestamps[$i]
}

# Add the new "metadata" worksheet as an additional tab (after "data").
$wsMeta = $wb.Worksheets.Add($null, $wsData)
$wsMeta.Name = "metadata"

# Header row
$wsMeta.Cells.Item(1, 2).Value = "data_name"
$wsMeta.Cells.Item(1, 3).Value = "data_id"
$wsMeta.Cells.Item(1, 4).Value = "data_version"
$wsMeta.Cells.Item(1, 5).Value = "data_version_created"
$wsMeta.Cells.Item(1, 6).Value = "panel_query_time"
$wsMeta.Cells.Item(1, 7).Value = "panel_get_request"

# Match the header formatting (bold font + border + centered/top alignment)
# used by the "data" sheet's header row.
$wsData.Range("B1").Copy()
[void]$wsMeta.Range("B1:G1").PasteSpecial(-4122)

# Data row
$wsMeta.Cells.Item(2, 1).Value = 0

# Match the formatting used by the "data" sheet's index column.
$wsData.Range("A2").Copy()
[void]$wsMeta.Range("A2").PasteSpecial(-4122)

$wsMeta.Cells.Item(2, 2).Value = "Autosomal recessive congenital ichthyosis"
$wsMeta.Cells.Item(2, 3).Value = 282

# "1.13" must stay a text value (not be coerced to the number 1.13), but
# without leaving any cell-level formatting behind.
$wsMeta.Cells.Item(2, 4).Value = "'1.13"
$wsMeta.Cells.Item(2, 4).Style = "Normal"

$wsMeta.Cells.Item(2, 5).Value = "2021-08-24T15:53:02.191596Z"
$wsMeta.Cells.Item(2, 6).Value = "2021-10-05 14:19:16.017292"
$wsMeta.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/282/?format=json"

[void]$wsData.Activate()
[void]$wsData.Range("A1").Select()
